$wb = $excel.ActiveWorkbook

# Update "想去人数" (interest count) figures that changed between scrapes.
# Both the "展览" sheet and the aggregated "全部类型" sheet contain the
# same rows and need the same update.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 3261
    $ws.Range("F5").Value = 1182
}
